$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new leader "Admin" in column A of row 3
$ws.Range("A3").Value = "Admin"

# Update the selection to match the authored state (active cell B3)
$ws.Range("B3").Select()
